# updated GSC export files
#
# The GSC HTTPS export gained one more day of data (2025-11-09). Append it
# as a new row at the bottom of the "Chart" sheet, right after the last
# existing date row, following the same Date | Non-HTTPS URLs | HTTPS URLs
# layout as every row above it.
#
# The "Table" sheet's header ("Issue" / "Validation" / "Pages") doesn't
# change content-wise -- in the source diff its shared-string indices only
# shift because the new date string gets interned into the shared-strings
# table, not because anything on that sheet was edited. So it needs no
# explicit edit here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

$colA = $ws.Cells.Item($newRow, 1)
$colB = $ws.Cells.Item($newRow, 2)
$colC = $ws.Cells.Item($newRow, 3)

# Column A holds the date as literal text (matches the existing A2:A34
# cells, which are shared strings, not real dates). A leading apostrophe
# makes Excel store the value as text instead of auto-converting it to a
# date serial; ClearFormats() then drops the resulting "quote prefix" cell
# style so the new cell keeps the workbook's plain default look, same as
# every other row.
$colA.Value = "'2025-11-09"
$colA.ClearFormats()

$colB.Value = 0
$colC.Value = 83
